$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2287.8333
$ws.Range("J19").Value = 664.8570999999999
$ws.Range("L19").Value = 664.8570999999999
$ws.Range("N19").Value = -1014.8571
$ws.Range("H33").Value = 136.73334
$ws.Range("I33").Value = 142.92857
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 142.92857
$ws.Range("L33").Value = 50
$ws.Range("M33").Value = 86.07142999999999
$ws.Range("N33").Value = -508
$ws.Range("H106").Value = 17546048
$ws.Range("I106").Value = 23811180
$ws.Range("K106").Value = 23811180
$ws.Range("M106").Value = -23810549
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""
$ws.Range("H132").Value = 3238.6667
$ws.Range("I132").Value = 3341.76
$ws.Range("J132").Value = 1950
$ws.Range("K132").Value = 10025.28
$ws.Range("L132").Value = 5850
$ws.Range("M132").Value = -7495.280000000001
$ws.Range("N132").Value = -10910
$ws.Range("H138").Value = 2220.1082
$ws.Range("I138").Value = 2126.3572
$ws.Range("J138").Value = 2241.9834
$ws.Range("K138").Value = 6379.071599999999
$ws.Range("L138").Value = 6725.9502
$ws.Range("M138").Value = -1239.071599999999
$ws.Range("N138").Value = -17005.9502

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5240.6924
$ws.Range("I32").Value = 5186.5884
$ws.Range("K32").Value = 5186.5884
$ws.Range("M32").Value = -4899.5884
$ws.Range("H61").Value = 1195.6389
$ws.Range("I61").Value = 1101.2285
$ws.Range("J61").Value = 4500
$ws.Range("K61").Value = 1101.2285
$ws.Range("L61").Value = 4500
$ws.Range("M61").Value = -889.2284999999999
$ws.Range("N61").Value = -4924
$ws.Range("H74").Value = 66667624
$ws.Range("I74").Value = 76924024
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 76924024
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = -76923150
$ws.Range("N74").Value = -2748
$ws.Range("H77").Value = 66667624
$ws.Range("I77").Value = 76924024
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 384620120
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = -384615752
$ws.Range("N77").Value = -13736
$ws.Range("H132").Value = 25487.92
$ws.Range("I132").Value = 1419.7307
$ws.Range("K132").Value = 4259.1921
$ws.Range("M132").Value = -1729.1921
$ws.Range("H135").Value = 53514.668
$ws.Range("J135").Value = 53514.668
$ws.Range("L135").Value = 53514.668
$ws.Range("N135").Value = -63654.668
$ws.Range("H136").Value = 1195.6389
$ws.Range("I136").Value = 1101.2285
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 3303.6855
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -753.6854999999996
$ws.Range("N136").Value = -18600

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 393.34784
$ws.Range("I94").Value = 411.8095
$ws.Range("J94").Value = 199.5
$ws.Range("K94").Value = 411.8095
$ws.Range("L94").Value = 199.5
$ws.Range("M94").Value = 39.19049999999999
$ws.Range("N94").Value = -1101.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 412.22223
$ws.Range("I22").Value = 338.75
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 338.75
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = 11.25
$ws.Range("N22").Value = -1700
$ws.Range("H31").Value = 15118.774
$ws.Range("I31").Value = 18501.348
$ws.Range("K31").Value = 18501.348
$ws.Range("M31").Value = -18206.348
$ws.Range("H34").Value = 15118.774
$ws.Range("I34").Value = 18501.348
$ws.Range("K34").Value = 18501.348
$ws.Range("M34").Value = -18299.348
$ws.Range("H132").Value = 10505.25
$ws.Range("I132").Value = 12486.363
$ws.Range("J132").Value = 3241.1667
$ws.Range("K132").Value = 37459.089
$ws.Range("L132").Value = 9723.500100000001
$ws.Range("M132").Value = -34929.089
$ws.Range("N132").Value = -14783.5001
$ws.Range("H134").Value = 658.95
$ws.Range("I134").Value = 535.4194
$ws.Range("K134").Value = 1606.2582
$ws.Range("M134").Value = 928.7418

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 68.64286
$ws.Range("I12").Value = 18
$ws.Range("J12").Value = 88.90000000000001
$ws.Range("K12").Value = 54
$ws.Range("L12").Value = 266.7
$ws.Range("M12").Value = 119
$ws.Range("N12").Value = -612.7
$ws.Range("H123").Value = 3866.5
$ws.Range("I123").Value = 2154.5
$ws.Range("J123").Value = 4722.5
$ws.Range("K123").Value = 6463.5
$ws.Range("L123").Value = 14167.5
$ws.Range("M123").Value = -4013.5
$ws.Range("N123").Value = -19067.5
$ws.Range("H131").Value = 159557.55
$ws.Range("J131").Value = 167507.75
$ws.Range("L131").Value = 502523.25
$ws.Range("N131").Value = -512603.25
$ws.Range("H132").Value = 875.1539
$ws.Range("I132").Value = 499.66666
$ws.Range("J132").Value = 1197
$ws.Range("K132").Value = 4496.99994
$ws.Range("L132").Value = 10773
$ws.Range("M132").Value = -1966.99994
$ws.Range("N132").Value = -15833

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20166.666
$ws.Range("I70").Value = 50500
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 50500
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -50230
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 20166.666
$ws.Range("I73").Value = 50500
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 50500
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -49564
$ws.Range("N73").Value = -6872
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = ""
$ws.Range("H107").Value = 5494754
$ws.Range("I107").Value = 283.7
$ws.Range("J107").Value = 19230930
$ws.Range("K107").Value = 283.7
$ws.Range("L107").Value = 19230930
$ws.Range("M107").Value = 1636.3
$ws.Range("N107").Value = -19234770
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1333.4839
$ws.Range("I136").Value = 1101.5264
$ws.Range("J136").Value = 1700.75
$ws.Range("K136").Value = 3304.5792
$ws.Range("L136").Value = 5102.25
$ws.Range("M136").Value = -754.5792000000001
$ws.Range("N136").Value = -10202.25

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1232.3334
$ws.Range("I96").Value = 800
$ws.Range("J96").Value = 1448.5
$ws.Range("K96").Value = 800
$ws.Range("L96").Value = 1448.5
$ws.Range("M96").Value = 573
$ws.Range("N96").Value = -4194.5
$ws.Range("H126").Value = 1882.2667
$ws.Range("I126").Value = 1144.1666
$ws.Range("J126").Value = 4834.6665
$ws.Range("K126").Value = 3432.4998
$ws.Range("L126").Value = 14503.9995
$ws.Range("M126").Value = -962.4998000000001
$ws.Range("N126").Value = -19443.9995
$ws.Range("H136").Value = 21740328
$ws.Range("I136").Value = 24391272
$ws.Range("J136").Value = 2580
$ws.Range("K136").Value = 73173816
$ws.Range("L136").Value = 7740
$ws.Range("M136").Value = -73171266
$ws.Range("N136").Value = -12840
